$p = $ppt.ActivePresentation

# --- Slide 1: subtitle "Bài " + "19. " (two runs) -> merge into a single run "Bài 19. " ---
$s1 = $p.Slides.Item(1)
$shp1 = $s1.Shapes.Item(2)
$tr1 = $shp1.TextFrame.TextRange
$tr1.Characters(1, 8).Text = "Bài 19. "

# --- Slide 17: title "Bài tập" (single run) -> split into "Bài " + "tập 19.1" ---
$s17 = $p.Slides.Item(17)
$shp17 = $s17.Shapes.Item(1)
$tr17 = $shp17.TextFrame.TextRange
$tr17.Characters(5, 3).Text = "tập 19.1"

# --- Slide 18: title "Bài tập (2)" (single run) -> split into "Bài tập " + "19.2" ---
$s18 = $p.Slides.Item(18)
$shp18 = $s18.Shapes.Item(1)
$tr18 = $shp18.TextFrame.TextRange
$tr18.Characters(9, 3).Text = "19.2"
